# Apply the SWAWR row reordering / renumbering edit to sheet1 (A10:F17).
# Rows 10-17 (the D="5" / SWAWR group) are reordered and the running
# sequence number in column E is renumbered 1-8 across the group.
#
# Values are written via a text-producing formula (="...") and then
# "frozen" in place with Copy + PasteSpecial(xlPasteValues). This keeps
# numeric-looking strings (e.g. "20101460", "5") stored as genuine text
# (shared-string) cells without Excel silently re-interpreting them as
# numbers, and - unlike toggling NumberFormat to "@" - it does not leave
# behind any extra/unused cell style entries, so existing cell styles
# (s="1") are preserved exactly as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# Final values for rows 10-17, columns A-F.
$rows = @(
    @{ A = "20101460"; B = "CHOMP ICE.CRM STR 50"; C = "SWAWR"; D = "5"; E = "1"; F = "RT,(E-1B)" },
    @{ A = "20069650"; B = "CHOMP ICE.CRM BLU 50"; C = "SWAWR"; D = "5"; E = "2"; F = "RT,(E-2B)" },
    @{ A = "20096157"; B = "CHOMP MLLOW STRW 50"; C = "SWAWR"; D = "5"; E = "3"; F = "RT,(E-1B)" },
    @{ A = "20078838"; B = "CHOMP GOLD.COIN 6'S"; C = "SWAWR"; D = "5"; E = "4"; F = "RT,(E-2B)" },
    @{ A = "20064556"; B = "CHOMP2 TWISTER 50G"; C = "SWAWR"; D = "5"; E = "5"; F = "RT,(E-1B)" },
    @{ A = "20137587"; B = "CHOMP2 MALLOWPOP 16G"; C = "SWAWR"; D = "5"; E = "6"; F = "RT,(E-1B)" },
    @{ A = "20133216"; B = "CHOMP TRIO PLAIN 21G"; C = "SWAWR"; D = "5"; E = "7"; F = "RT,(E-1B)" },
    @{ A = "20133246"; B = "CHOMP TRIO RNBOW 21G"; C = "SWAWR"; D = "5"; E = "8"; F = "RT,(E-1B)" }
)

$startRow = 10
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    Set-TextValue $ws.Cells.Item($r, 1) $row.A
    Set-TextValue $ws.Cells.Item($r, 2) $row.B
    Set-TextValue $ws.Cells.Item($r, 3) $row.C
    Set-TextValue $ws.Cells.Item($r, 4) $row.D
    Set-TextValue $ws.Cells.Item($r, 5) $row.E
    Set-TextValue $ws.Cells.Item($r, 6) $row.F
}
